$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "checklist"

# New column G holds the product code; format as text so values like
# "6-1" are not coerced into dates. Header (G1) keeps the default/bold
# header style, so only format the data rows.
$ws.Range("G2:G35").NumberFormat = "@"

$ws.Range("G1").Value = "product_code"

$ws.Range("G2").Value = "6-1"
$ws.Range("G3").Value = "6-5"
$ws.Range("G4").Value = "6-5"
$ws.Range("G5").Value = "6-28"
$ws.Range("G7").Value = "6-35"
$ws.Range("G12").Value = "6-2"
$ws.Range("G13").Value = "6-3"
$ws.Range("G14").Value = "6-8"
$ws.Range("G15").Value = "6-9"
$ws.Range("G16").Value = "6-11"
$ws.Range("G17").Value = "6-13"
$ws.Range("G18").Value = "6-15"
$ws.Range("G20").Value = "6-20"
$ws.Range("G21").Value = "6-27"
$ws.Range("G23").Value = "6-4"
$ws.Range("G24").Value = "6-7"
$ws.Range("G25").Value = "6-10"

# Match the author's final selection position
[void]$ws.Range("G22").Select()
